# The deck's slide master currently uses the "Integral" theme
# (ppt/theme/theme2.xml) while ppt/theme/theme1.xml (used only by the
# notes master) holds the stock "Office Theme" colours. The authored
# edit swaps the two themes' contents, so the slides end up using the
# "Office Theme" colour palette.
#
# The PowerPoint object model only exposes one editable Theme (the one
# tied to the slide master / presentation theme relationship), reached
# through Master.Theme / Slide.Master.Theme / Design.SlideMaster.Theme.
# Its ThemeColorScheme collection holds the 12 theme colours in the
# fixed order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink. Re-point
# each of those RGB values at the stock "Office Theme" palette to apply
# the swap.

$p = $ppt.ActivePresentation
$master = $p.Slides.Item(1).Master
$theme = $master.Theme
$colorScheme = $theme.ThemeColorScheme

# [index] = decimal RGB() value (0x00BBGGRR, PowerPoint's native colour
# encoding) for the "Office Theme" colour scheme, in dk1, lt1, dk2, lt2,
# accent1, accent2, accent3, accent4, accent5, accent6, hlink, folHlink
# order.
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeColors[$i - 1]
}
